# RECIBOS NOMINA - update week header from SEMANA 35 (23-29 Agosto 2021)
# to SEMANA 37 (06-12 Septiembre 2021), bump extra-hours count/amount,
# bump the "media" (medio dia) allowance, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Week header text (B9 is the single source; H9/B27/H27/B43/H43/B60 all
# reference it via formulas and pick up the change automatically).
$ws.Range("B9").Value = "SEMANA   37  DEL  06  Al    12 DE SEPTIEMBRE          2021"

# Extra hours: 6 -> 8, and the amount becomes a formula (280 per hour)
# instead of a hard-coded number.
$ws.Range("J21").Value = 8
$ws.Range("K21").Formula = "=280*J21"

# "Media" (medio dia) allowance bump for the second pay period.
$ws.Range("K40").Value = 1250

# Move the active selection / view to match where the editor left off.
$null = $ws.Range("K41").Select()
